# Update the "取得日時" (retrieved datetime) column for existing rows 2-20
# on the "ランサーズ" sheet to the new timestamp, as if a fresh scrape
# appended at 2025-10-03 12:44:32 JST overwrote the previous batch's
# timestamp value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-03 12:44:32"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 20 }

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
